# Update the Product Register script
# - Update the From/To date sample values
# - Rename the "Invoice.number.heading" key to "Invoice.name.heading"
# - Append new automation steps (rows 104-111) for Godown entry,
#   filter table, invoice buyer address and search value checks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values -------------------------------------------------
$ws.Range("B74").Value() = "15-11-2023"
$ws.Range("B77").Value() = "20-11-2023"
$ws.Range("A102").Value() = "Invoice.name.heading"

# --- Append new rows ---------------------------------------------------------
$newRows = @(
  @("Godown.enter.values", "Primary Godown (Primary Godown's address)"),
  @("godown.click.first.values", "//main[@class='mb-5']/section/div/div/div[2]/form/fieldset[2]/div/div/div/div/button[1]"),
  @("filter.Table", "//div[@class='popover-body']"),
  @("filter.Table.heading", "//div[@class='popover-body']/div[text()=' Invoice Type ']"),
  @("Invoice.buyer.address", "//main[@class='mb-5']/section/div/div[2]/div/div[2]/div/div/div/div[3]/div/div/p"),
  @("searchValues", "1/SL-23")
)

$r = 104
foreach ($pair in $newRows) {
  $ws.Cells.Item($r, 1).Value() = $pair[0]
  $ws.Cells.Item($r, 2).Value() = $pair[1]
  $r = $r + 1
}

# Row 110 only has a key in column A (no value in column B)
$ws.Cells.Item(110, 1).Value() = "search.invalid.values"

$ws.Cells.Item(111, 1).Value() = "filter.Popup.Message"
$ws.Cells.Item(111, 2).Value() = "//div[@class='toast-body'][contains(text(),'No data present for selected filter')]"

# --- Update the view so the new rows are visible -----------------------------
$ws.Range("B98").Select()
$win = $excel.ActiveWindow
$win.ScrollRow() = 84
